# Update cryptos list with latest prices and volume percentages
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    # Force the cell to remain a text value (matching the source data,
    # which stores these numeric-looking price strings as text) instead
    # of letting Excel auto-convert it to a number.
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value2 = $val
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue 2 4 '38.762.85'
$ws.Cells.Item(2, 5).Value = '  +2.76%  '

# Row 3
Set-TextValue 3 4 '2.091.10'
$ws.Cells.Item(3, 5).Value = '  +2.57%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.01%  '

# Row 5
Set-TextValue 5 4 '228.03'
$ws.Cells.Item(5, 5).Value = '  +0.32%  '

# Row 6
$ws.Cells.Item(6, 5).Value = '  +0.86%  '

# Row 7
Set-TextValue 7 4 '60.39'
$ws.Cells.Item(7, 5).Value = '  +1.32%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  -0.01%  '

# Row 9
$ws.Cells.Item(9, 5).Value = '  +2.02%  '

# Row 10
Set-TextValue 10 4 '0.0831'
$ws.Cells.Item(10, 5).Value = '  -1.05%  '

# Row 11
$ws.Cells.Item(11, 5).Value = '  -0.17%  '

# Row 12
Set-TextValue 12 4 '2.402.61'
$ws.Cells.Item(12, 5).Value = '  +2.69%  '

# Row 13
Set-TextValue 13 4 '14.92'
$ws.Cells.Item(13, 5).Value = '  +3.26%  '

# Row 14
Set-TextValue 14 4 '22.05'
$ws.Cells.Item(14, 5).Value = '  +4.83%  '

# Row 15
Set-TextValue 15 4 '0.796'
$ws.Cells.Item(15, 5).Value = '  +3.18%  '

# Row 16
$ws.Cells.Item(16, 5).Value = '  -0.08%  '

# Row 17
Set-TextValue 17 4 '2.124.72'
$ws.Cells.Item(17, 5).Value = '  +3.88%  '

# Row 18
Set-TextValue 18 4 '38.713.23'
$ws.Cells.Item(18, 5).Value = '  +2.65%  '

# Row 19
Set-TextValue 19 4 '71.94'
$ws.Cells.Item(19, 5).Value = '  +3.56%  '

# Row 20
Set-TextValue 20 4 '6.04'
$ws.Cells.Item(20, 5).Value = '  +2.24%  '

# Row 21
Set-TextValue 21 4 '0.0₃0834'
$ws.Cells.Item(21, 5).Value = '  +1.31%  '

# Row 22
Set-TextValue 22 4 '226.10'
$ws.Cells.Item(22, 5).Value = '  +1.04%  '

# Row 23
Set-TextValue 23 4 '0.999'
$ws.Cells.Item(23, 5).Value = '  -0.54%  '

# Row 24
Set-TextValue 24 4 '2.38'
$ws.Cells.Item(24, 5).Value = '  -0.38%  '

# Row 25
Set-TextValue 25 4 '2.32'
$ws.Cells.Item(25, 5).Value = '  +2.04%  '

# Row 26
Set-TextValue 26 4 '170.47'
$ws.Cells.Item(26, 5).Value = '  +1.26%  '

# Row 27
$ws.Cells.Item(27, 5).Value = '  +1.24%  '

# Row 28
$ws.Cells.Item(28, 5).Value = '  +6.07%  '

# Row 29
$ws.Cells.Item(29, 5).Value = '  +9.25%  '

# Row 30
$ws.Cells.Item(30, 5).Value = '  +1.85%  '

# Row 31
$ws.Cells.Item(31, 5).Value = '  +0.90%  '

# Row 32
$ws.Cells.Item(32, 5).Value = '  +4.89%  '

# Row 33
$ws.Cells.Item(33, 5).Value = '  +5.33%  '

# Row 34
Set-TextValue 34 4 '4.47'
$ws.Cells.Item(34, 5).Value = '  +2.07%  '

# Row 35
$ws.Cells.Item(35, 5).Value = '  +1.22%  '

# Row 36
Set-TextValue 36 4 '2.38'
$ws.Cells.Item(36, 5).Value = '  +2.50%  '

# Row 37
Set-TextValue 37 4 '6.40'
$ws.Cells.Item(37, 5).Value = '  +0.11%  '

# Row 38
$ws.Cells.Item(38, 5).Value = '  +2.78%  '

# Row 39
$ws.Cells.Item(39, 5).Value = '  +0.00%  '

# Row 40
$ws.Cells.Item(40, 5).Value = '  +1.45%  '

# Row 41
$ws.Cells.Item(41, 2).Value = 'Maker'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue 41 4 '1.537.60'
$ws.Cells.Item(41, 5).Value = '  +0.61%  '

# Row 42
$ws.Cells.Item(42, 2).Value = 'Aave'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue 42 4 '100.97'
$ws.Cells.Item(42, 5).Value = '  +3.80%  '

# Row 44
Set-TextValue 44 4 '0.0926'
$ws.Cells.Item(44, 5).Value = '  +2.26%  '

# Row 45
$ws.Cells.Item(45, 5).Value = '  -1.04%  '

# Row 46
Set-TextValue 46 4 '7.65'
$ws.Cells.Item(46, 5).Value = '  +9.40%  '

# Row 47
Set-TextValue 47 4 '1.11'
$ws.Cells.Item(47, 5).Value = '  +0.51%  '

# Row 48
Set-TextValue 48 4 '4.11'
$ws.Cells.Item(48, 5).Value = '  -3.67%  '

# Row 49
$ws.Cells.Item(49, 5).Value = '  +2.80%  '

# Row 50
Set-TextValue 50 4 '2.97'
$ws.Cells.Item(50, 5).Value = '  +1.09%  '

# Row 51
Set-TextValue 51 4 '2.288.88'
$ws.Cells.Item(51, 5).Value = '  +2.69%  '

